# Add 2022-Q3 data
#
# 1. "总计" (summary) sheet: insert a new row for "2022-Q3" right after the
#    header row, shifting the existing quarters down, and renumber the
#    leading index column.
# 2. Insert a brand-new worksheet named "2022-Q3" right after "总计" (and
#    before "2022-Q2") holding the per-fund detail rows for that quarter.

$wb = $excel.ActiveWorkbook

# --- 1. Update the "总计" overview sheet -----------------------------------
$summary = $wb.Worksheets.Item(1)

$summary.Rows.Item(2).Insert()
$summary.Range("A3:D3").Copy($summary.Range("A2:D2"))

$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = "2022-Q3"
$summary.Cells.Item(2, 3).Value = 2
$summary.Cells.Item(2, 4).Value = 0.01

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(5, 1).Value = 3

# --- 2. Insert the new "2022-Q3" detail sheet ------------------------------
$afterSheet = $wb.Worksheets.Item(1)
$q3 = $wb.Worksheets.Add($null, $afterSheet)
$q3.Name = "2022-Q3"

$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"
$q3.Range("B1:H1").Font.Bold = $true
$q3.Range("B1:H1").HorizontalAlignment = -4108
$q3.Range("B1:H1").VerticalAlignment = -4160

$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "003242"
$q3.Cells.Item(2, 3).Value = "创金合信量化发现灵活配置混合C"
$q3.Cells.Item(2, 4).Value = "0.40"
$q3.Cells.Item(2, 5).Value = "92.08"
$q3.Cells.Item(2, 6).Value = "1.44"
$q3.Cells.Item(2, 7).Value = "0.0058"
$q3.Cells.Item(2, 8).Value = 10
$q3.Cells.Item(2, 1).Font.Bold = $true
$q3.Cells.Item(2, 1).HorizontalAlignment = -4108
$q3.Cells.Item(2, 1).VerticalAlignment = -4160

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "003241"
$q3.Cells.Item(3, 3).Value = "创金合信量化发现灵活配置混合A"
$q3.Cells.Item(3, 4).Value = "0.32"
$q3.Cells.Item(3, 5).Value = "92.08"
$q3.Cells.Item(3, 6).Value = "1.44"
$q3.Cells.Item(3, 7).Value = "0.0046"
$q3.Cells.Item(3, 8).Value = 10
$q3.Cells.Item(3, 1).Font.Bold = $true
$q3.Cells.Item(3, 1).HorizontalAlignment = -4108
$q3.Cells.Item(3, 1).VerticalAlignment = -4160

$q3.Range("A1").Select()
